$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '43.708.62'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '2.287.21'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('E4').Value = '  +0.46%  '
Set-TextValue $ws 'D5' '111.11'
$ws.Range('E5').Value = '  +15.16%  '
Set-TextValue $ws 'D6' '267.33'
$ws.Range('E6').Value = '  -0.30%  '
Set-TextValue $ws 'D7' '0.624'
$ws.Range('E7').Value = '  +1.19%  '
Set-TextValue $ws 'D9' '0.619'
$ws.Range('E9').Value = '  +1.26%  '
Set-TextValue $ws 'D10' '47.50'
$ws.Range('E10').Value = '  +3.61%  '
Set-TextValue $ws 'D11' '0.0943'
$ws.Range('E11').Value = '  +0.80%  '
Set-TextValue $ws 'D12' '8.84'
$ws.Range('E12').Value = '  +11.33%  '
Set-TextValue $ws 'D13' '0.106'
$ws.Range('E13').Value = '  +1.20%  '
Set-TextValue $ws 'D14' '15.74'
$ws.Range('E14').Value = '  +2.60%  '
$ws.Range('D15').Value = '2.633.86'
$ws.Range('E15').Value = '  -0.08%  '
Set-TextValue $ws 'D16' '0.844'
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('D17').Value = '2.296.66'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').Value = '43.727.96'
$ws.Range('E19').Value = '  +1.51%  '
Set-TextValue $ws 'D20' '6.59'
$ws.Range('E20').Value = '  +6.21%  '
Set-TextValue $ws 'D21' '72.34'
$ws.Range('E21').Value = '  +0.01%  '
Set-TextValue $ws 'D22' '2.45'
$ws.Range('E22').Value = '  -3.96%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws 'D23' '232.16'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws 'D24' '9.61'
$ws.Range('E24').Value = '  +5.63%  '
Set-TextValue $ws 'D25' '2.79'
$ws.Range('E25').Value = '  +9.96%  '
Set-TextValue $ws 'D26' '0.999'
$ws.Range('E26').Value = '  -0.04%  '
Set-TextValue $ws 'D27' '11.60'
$ws.Range('E27').Value = '  +3.44%  '
Set-TextValue $ws 'D28' '41.59'
$ws.Range('E28').Value = '  +3.74%  '
Set-TextValue $ws 'D29' '3.41'
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('E30').Value = '  +1.70%  '
Set-TextValue $ws 'D31' '176.22'
$ws.Range('E31').Value = '  +0.33%  '
Set-TextValue $ws 'D32' '0.0928'
$ws.Range('E32').Value = '  +3.33%  '
Set-TextValue $ws 'D33' '21.53'
$ws.Range('E33').Value = '  -2.30%  '
Set-TextValue $ws 'D34' '5.62'
$ws.Range('E34').Value = '  +4.74%  '
$ws.Range('E35').Value = '  +0.87%  '
Set-TextValue $ws 'D36' '4.71'
$ws.Range('E36').Value = '  +7.76%  '
Set-TextValue $ws 'D37' '0.0361'
$ws.Range('E37').Value = '  +2.04%  '
$ws.Range('E38').Value = '  -0.52%  '
Set-TextValue $ws 'D39' '3.75'
$ws.Range('E39').Value = '  +11.03%  '
Set-TextValue $ws 'D40' '0.242'
$ws.Range('E40').Value = '  -0.16%  '
Set-TextValue $ws 'D41' '13.74'
$ws.Range('E41').Value = '  +11.50%  '
Set-TextValue $ws 'D42' '2.39'
$ws.Range('E42').Value = '  +3.46%  '
Set-TextValue $ws 'D43' '71.14'
$ws.Range('E43').Value = '  +9.29%  '
Set-TextValue $ws 'D44' '6.15'
$ws.Range('E44').Value = '  +18.26%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws 'D45' '1.00'
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws 'D46' '1.40'
$ws.Range('E46').Value = '  +3.08%  '
Set-TextValue $ws 'D47' '8.86'
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('E48').Value = '  -1.08%  '
Set-TextValue $ws 'D49' '101.44'
$ws.Range('E49').Value = '  +4.23%  '
Set-TextValue $ws 'D50' '1.22'
$ws.Range('E50').Value = '  +2.56%  '
Set-TextValue $ws 'D51' '0.446'
$ws.Range('E51').Value = '  +6.41%  '
